$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.430.03"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.867.58"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "BNB"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "243.35"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.44%  "

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "XRP"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7058"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.70%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07861"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.78%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3128"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.12%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "24.49"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07996"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.895.61"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.203"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.93%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "93.34"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.26%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6995"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.17%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.481"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.66%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "29.509.02"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008378"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.34"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.67%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.139.58"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.49%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.10"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.42%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.04%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.627"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.32%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.001"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.022"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.55%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.09"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -1.00%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.98%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.500"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.21%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.323"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.26%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.266"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.205"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +1.10%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05307"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.61%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7494"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.82%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.171"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.01%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +0.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01882"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.21%  "

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.75%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.748"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8957"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.109"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.77%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "109.01"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.68%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "71.40"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.72%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.12%  "

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -3.47%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.036.19"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.33%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.570"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.21%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.791"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.5174"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.01%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4304"
